$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data
$ws.Range("D2").Value = '29.319.03'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.874.79'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'0.7099"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = "'242.19"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.07791"
$ws.Range("E8").Value = '  +0.70%  '
$ws.Range("D9").Value = "'0.3109"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = "'25.12"
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("D11").Value = "'0.08422"
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").Value = '1.864.27'
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").Value = "'5.244"
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").Value = "'0.7182"
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = "'91.17"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = '29.320.61'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = "'6.094"
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("D18").Value = "'0.000008292"
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").Value = "'240.84"
$ws.Range("E19").Value = '  -1.35%  '
$ws.Range("D20").Value = "'13.24"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '2.121.24'
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = "'7.754"
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = "'9.042"
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'162.24"
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").Value = "'18.52"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = "'1.508"
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = "'4.411"
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'4.327"
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = "'1.288"
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = "'0.05388"
$ws.Range("E33").Value = '  +3.62%  '
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = "'0.7522"
$ws.Range("E35").Value = '  -3.02%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = "'1.178"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").Value = "'2.685"
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = "'0.01890"
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").Value = '1.237.40'
$ws.Range("E39").Value = '  +6.89%  '
$ws.Range("D40").Value = "'2.737"
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("D41").Value = "'6.487"
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("D42").Value = "'0.8940"
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  -1.30%  '
$ws.Range("D44").Value = "'108.98"
$ws.Range("E44").Value = '  +5.01%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '2.019.35'
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("E47").Value = '  +9.51%  '
$ws.Range("D48").Value = "'1.800"
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").Value = "'0.5200"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = "'9.463"
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("D51").Value = "'0.4344"
$ws.Range("E51").Value = '  +0.76%  '

# Restore default (General) style on cells forced to text so only the value differs
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
